$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: remove any hyperlink anchored on a given row of a worksheet, then
# delete the entire row (shifting everything below it up).
# ---------------------------------------------------------------------------
function Remove-RowAndHyperlinks {
    param($ws, [int]$row)

    # Address looks like "$A$3" - match the trailing "$<row>" part. Deleting
    # from the live collection while a `foreach` enumerator is positioned on
    # it skips entries, so rescan-and-delete-first-match instead.
    $suffix = '$' + $row
    $more = $true
    while ($more) {
        $more = $false
        foreach ($hl in $ws.Hyperlinks) {
            $addr = $hl.Range.Address()
            if ($addr.EndsWith($suffix)) {
                $hl.Delete()
                $more = $true
                break
            }
        }
    }

    $ws.Rows.Item($row).Delete()
}

# ---------------------------------------------------------------------------
# Overview sheet: drop the f42e4ba0... row (row 3).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Remove-RowAndHyperlinks $wsOverview 3

# ---------------------------------------------------------------------------
# zh-cn sheet: refresh the handoff/handback timestamps on row 2, then drop
# the f42e4ba0... row (row 3).
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Cells.Item(2, 5).Value = "2016-03-30 10:11:46"
$wsZh.Cells.Item(2, 8).Value = "2016-03-30 10:12:45"
Remove-RowAndHyperlinks $wsZh 3

# ---------------------------------------------------------------------------
# de-de sheet: refresh the handoff/handback timestamps on row 2, then drop
# the f42e4ba0... row (row 3).
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Cells.Item(2, 5).Value = "2016-03-30 10:11:57"
$wsDe.Cells.Item(2, 8).Value = "2016-03-30 10:13:04"
Remove-RowAndHyperlinks $wsDe 3
